$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.608.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.74%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.314.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.53%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "181.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -10.69%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "531.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.32%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.606"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.57%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.309.89"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.37%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.621"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.38%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.99"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.28%  "

# Row 12
$ws.Range("E12").Value = "  -7.13%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000265"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.18%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.826.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.91%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.302.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.99%  "

# Row 17
$ws.Range("E17").Value = "  -5.63%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.72%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "64.327.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.69%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.93%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.963"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.52%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "375.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.11%  "

# Row 23
$ws.Range("E23").Value = "  -5.71%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.40%  "

# Row 25
$ws.Range("E25").Value = "  -6.80%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.83%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.32%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.20%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.09%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.49%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "29.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.69%  "

# Row 32
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "643.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.13%  "

# Row 33
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.76%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.59%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.79%  "

# Row 36
$ws.Range("E36").Value = "  -6.29%  "

# Row 37
$ws.Range("E37").Value = "  +0.05%  "

# Row 38
$ws.Range("E38").Value = "  -3.25%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.54%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0739"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.49%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.997"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "

# Row 42
$ws.Range("E42").Value = "  -3.70%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.907.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.59%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.07%  "

# Row 45
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.97%  "

# Row 46
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0404"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.32%  "

# Row 47
$ws.Range("E47").Value = "  +10.03%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.66%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.56%  "

# Row 50
$ws.Range("E50").Value = "  -1.06%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.28%  "
